# Remove the "syntok" dependency row from the Acknowledgments sheet.
# This corresponds to removing syntok's sentence segmenter / word tokenizer
# references from the Wordless acknowledgments workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the row whose "Name" column (column A) contains "syntok" and delete
# the entire row, shifting the rows below upward (matches Excel's native
# row-delete behaviour, which keeps per-row formatting/hyperlinks intact).
$found = $ws.Range("A1:A1000").Find("syntok", [Type]::Missing, [Type]::Missing, 1)
if ($found -ne $null) {
    $ws.Rows.Item($found.Row).Delete() | Out-Null
}
